# Refresh the crypto symbol list (coinranking.com scrape) in-place:
#   - rows 9-15 and 41-42 are re-sorted, so the Coin/Link/Price/Volume
#     columns are rewritten with the coin that now belongs at that row
#   - every other changed row just gets fresh Price / Volume(1h) figures
# All of column D/E values are stored as literal text in the source file
# (e.g. "256.37", "0.52%"), so we force a text number-format before
# writing them and restore the default "Normal" style afterwards —
# otherwise Excel would silently convert these into real numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-valued cells (coin names / URLs) ---
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'

# --- Numeric-looking cells that must remain literal text (prices / percentages) ---
# Force text number-format first so Excel does not coerce these into real numbers,
# then restore the Normal style so no stray style index is left on the cell.
$numTextCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "E16", "D17", "E17", "E18", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "D25", "E25", "D26", "E26", "E27", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "E45", "D48", "E48")
foreach ($cellRef in $numTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '256.39'
$ws.Range("E2").Value = '0.56%'
$ws.Range("D3").Value = '26.91'
$ws.Range("E3").Value = '-4.72%'
$ws.Range("D4").Value = '4.747'
$ws.Range("E4").Value = '-9.27%'
$ws.Range("D5").Value = '0.05928'
$ws.Range("E5").Value = '1.24%'
$ws.Range("D6").Value = '6.658'
$ws.Range("E6").Value = '-1.13%'
$ws.Range("D7").Value = '0.8672'
$ws.Range("E7").Value = '0.13%'
$ws.Range("D8").Value = '0.9394'
$ws.Range("E8").Value = '-4.63%'
$ws.Range("D9").Value = '0.01048'
$ws.Range("E9").Value = '1,625.21%'
$ws.Range("D10").Value = '0.1400'
$ws.Range("E10").Value = '-0.70%'
$ws.Range("D11").Value = '0.03753'
$ws.Range("E11").Value = '7.70%'
$ws.Range("D12").Value = '0.07104'
$ws.Range("E12").Value = '-0.87%'
$ws.Range("D13").Value = '0.03163'
$ws.Range("E13").Value = '-0.64%'
$ws.Range("D14").Value = '0.09246'
$ws.Range("E14").Value = '0.19%'
$ws.Range("D15").Value = '0.001546'
$ws.Range("E15").Value = '-0.05%'
$ws.Range("E16").Value = '3.10%'
$ws.Range("D17").Value = '3.487'
$ws.Range("E17").Value = '-0.35%'
$ws.Range("E18").Value = '-0.63%'
$ws.Range("E19").Value = '1.65%'
$ws.Range("D20").Value = '0.3154'
$ws.Range("E20").Value = '-0.81%'
$ws.Range("E21").Value = '0.24%'
$ws.Range("D22").Value = '3.807'
$ws.Range("E22").Value = '7.14%'
$ws.Range("D23").Value = '0.04203'
$ws.Range("E23").Value = '1.04%'
$ws.Range("D25").Value = '0.001220'
$ws.Range("E25").Value = '-0.64%'
$ws.Range("D26").Value = '0.004284'
$ws.Range("E26").Value = '-10.68%'
$ws.Range("E27").Value = '-0.10%'
$ws.Range("E40").Value = '0.27%'
$ws.Range("D41").Value = '0.006166'
$ws.Range("E41").Value = '5.00%'
$ws.Range("D42").Value = '0.1102'
$ws.Range("E42").Value = '0.03%'
$ws.Range("D43").Value = '0.002199'
$ws.Range("E43").Value = '-4.44%'
$ws.Range("D44").Value = '0.01141'
$ws.Range("E44").Value = '17.79%'
$ws.Range("E45").Value = '5.00%'
$ws.Range("D48").Value = '0.002385'
$ws.Range("E48").Value = '10.98%'

foreach ($cellRef in $numTextCells) {
    $ws.Range($cellRef).Style = "Normal"
}
